$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is "b.md" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-13 14:44:29"

# --- zh-cn sheet: row 3 is "b.md" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C3").Value = "Ready for handoff"
# "False" is auto-typed to a Boolean by plain assignment; force text with a
# leading quote, then drop the resulting quote-prefix style so the cell
# keeps its original (default) formatting.
$ws2.Range("F3").Value = "'False"
$ws2.Range("F3").Style = "Normal"
$ws2.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-13 14:44:21"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/0079857b52f6cf8fbf39bbb293cb7f1f426a7165/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/a925cc76e234ec4214349e261c37a063d4f10ece/e2e/b.md."
$ws2.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is "b.md" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("F3").Value = "'False"
$ws3.Range("F3").Style = "Normal"
$ws3.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-13 14:44:29"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/0079857b52f6cf8fbf39bbb293cb7f1f426a7165/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/a925cc76e234ec4214349e261c37a063d4f10ece/e2e/b.md."
$ws3.Columns.Item(16).ColumnWidth = 39.166666666666664
